# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet right before
# column N ("Late"), shifting the existing Late / heading(Date) / Outstanding
# columns one place to the right, and re-select the Repayment schedule tab.

$wb = $excel.ActiveWorkbook

$repayment = $wb.Worksheets.Item("Repayment schedule")
$transactions = $wb.Worksheets.Item("Transactions")

# Insert a blank column before column N (shifts N->O, O->P, P->Q)
$repayment.Columns("N").Insert()

# The newly inserted column picks up the width of the column to its left (M)
$repayment.Columns("N").ColumnWidth = $repayment.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and update its selection
$repayment.Activate() | Out-Null
$repayment.Range("I20").Select() | Out-Null
